$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17 — this shifts the existing rows 17-112 down to 18-113
# (carrying their formatting/styles along, same as Excel's native row insert),
# and extends the used range from A1:T112 to A1:T113.
$ws.Rows("17").Insert()

# Populate the newly inserted row 17 with the new weekly price-report record.
$ws.Range("A17").Value = 10
$ws.Range("B17").Value = "Vega Modelo de Temuco"
$ws.Range("C17").Value = "La Araucanía"
$ws.Range("D17").Value = 44901
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100101
$ws.Range("H17").Value = "Berries"
$ws.Range("I17").Value = 100101001
$ws.Range("J17").Value = "Arándano (blue)"
$ws.Range("K17").Value = "Sin especificar"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 400
$ws.Range("N17").Value = 2000
$ws.Range("O17").Value = 2200
$ws.Range("P17").Value = 2100
$ws.Range("Q17").Value = '$/envase 1 kilo'
$ws.Range("R17").Value = "Región del Maule"
$ws.Range("S17").Value = 2100
$ws.Range("T17").Value = 1
